# Fruta / hortaliza, semanal
# Insert a new weekly record at row 82 (pushing the existing rows 82:113 down
# to 83:114) for "Macroferia Regional de Talca" / Apio / Americana (o) / Primera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 82:113 down one row, creating a blank row 82
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with this week's data
$ws.Range("A82").Value = 5
$ws.Range("B82").Value = "Macroferia Regional de Talca"
$ws.Range("C82").Value = "Maule"
$ws.Range("D82").Value = 44466
$ws.Range("E82").Value = 7
$ws.Range("F82").Value = 100112017
$ws.Range("G82").Value = "Apio"
$ws.Range("H82").Value = "Americana (o)"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 500
$ws.Range("K82").Value = 8000
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = 8000
$ws.Range("N82").Value = "`$/docena de matas"
$ws.Range("O82").Value = "Provincia del Elquí"
$ws.Range("P82").Value = 1333
$ws.Range("Q82").Value = 6
$ws.Range("R82").Value = "Hortaliza"
